# Update "想去人数" (interest count) figures in F column across sheets,
# matching the regenerated gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1148
$ws1.Range("F9").Value = 289
$ws1.Range("F13").Value = 170
$ws1.Range("F14").Value = 3758
$ws1.Range("F15").Value = 23
$ws1.Range("F20").Value = 297

# Sheet "演出" (rId2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 287

# Sheet "全部类型" (rId4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1148
$ws4.Range("F14").Value = 289
$ws4.Range("F18").Value = 170
$ws4.Range("F19").Value = 3758
$ws4.Range("F20").Value = 23
$ws4.Range("F26").Value = 297
$ws4.Range("F32").Value = 287
